# HSE MNIS / Converted.xlsx - fixes in conversion script
#
# The source currency-adjusted market-cap (C) and total-assets (D) columns
# were actually already local-currency figures that still needed the FX
# conversion applied. This script:
#   1. Adds two new columns, F (CRNCY_ADJ_MKT_CAP_USD) and G (BS_TOT_ASSET_USD),
#      preserving the original (pre-conversion) C/D values there.
#   2. Overwrites C and D with the USD-converted figures (value * rate),
#      rounded to 2 decimal places, matching the corrected conversion script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Clone the formatting of the existing header cell (E1) onto the two new
# header cells, then set their text.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(1, 6).Value = "CRNCY_ADJ_MKT_CAP_USD"
$ws.Cells.Item(1, 7).Value = "BS_TOT_ASSET_USD"

# --- Data rows -----------------------------------------------------------
for ($row = 2; $row -le 101; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)
    $gCell = $ws.Cells.Item($row, 7)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    if ($cVal -eq "") {
        $fCell.Value = ""
    } else {
        $fCell.Value = $cVal
        $cCell.Value = [Math]::Round($cVal * $eVal, 2)
    }

    if ($dVal -eq "") {
        $gCell.Value = ""
    } else {
        $gCell.Value = $dVal
        $dCell.Value = [Math]::Round($dVal * $eVal, 2)
    }
}
